$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ip_address_list": append extra note lines to D5, then drop the
# "515_ZF Stara Boleslav" row (row 12) entirely.
# ---------------------------------------------------------------------------
$wsIp = $wb.Worksheets.Item("ip_address_list")

$d5 = $wsIp.Range("D5")
$d5.Value = $d5.Value2 + "`nf`nffg`nf`nf"

$wsIp.Rows(12).Delete()

# ---------------------------------------------------------------------------
# Sheet "ip_adress_fav_list": append the same extra note lines to D2, then
# drop the "474 B_Austin" row (row 3) so "527_Teijin" (old row 4) shifts up
# to become row 3.
# ---------------------------------------------------------------------------
$wsFav = $wb.Worksheets.Item("ip_adress_fav_list")

$d2 = $wsFav.Range("D2")
$d2.Value = $d2.Value2 + "`nf`nffg`nf`nf"

$wsFav.Rows(3).Delete()

# ---------------------------------------------------------------------------
# Sheet "disk_list": new note in F3, reworked note text in F4/F5, renamed /
# re-pointed "518_Valeo" row (row 5), then drop the "474_B Austin" row
# (row 6) entirely.
# ---------------------------------------------------------------------------
$wsDisk = $wb.Worksheets.Item("disk_list")

$wsDisk.Range("F3").Value = "f`nffffffffffffffffff"

$newNote = "Druha sít, ixonah`ndasf`ndfa`ndfa`ndfadfaafd`ndaf`ndfa`ndfa"
$wsDisk.Range("F4").Value = $newNote
$wsDisk.Range("F5").Value = $newNote

$wsDisk.Range("A5").Value = "518_Valeo II8"
$wsDisk.Range("C5").Value = "\\192.168.1.10\10_vision"

$wsDisk.Rows(6).Delete()

# ---------------------------------------------------------------------------
# Sheet "Settings": fix the default-interface setting from 4 to 0.
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 0
